$wb = $excel.ActiveWorkbook

# --- Debit ---
$ws = $wb.Worksheets.Item("Debit")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Fri Aug 29 12:39:12 EDT 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Fri Aug 29 12:39:17 EDT 2025"
$ws.Range("A4").Value = "Pass"
$ws.Range("B4").Value = "Fri Aug 29 12:39:22 EDT 2025"
$ws.Range("A5").Value = "Pass"
$ws.Range("B5").Value = "Fri Aug 29 12:39:26 EDT 2025"
$ws.Range("A6").Value = "Pass"
$ws.Range("B6").Value = "Fri Aug 29 12:39:31 EDT 2025"
$ws.Range("A7").Value = "Pass"
$ws.Range("B7").Value = "Fri Aug 29 12:39:36 EDT 2025"
$ws.Range("A8").Value = "Pass"
$ws.Range("B8").Value = "Fri Aug 29 12:39:41 EDT 2025"

# --- Debit-ZeroDollar ---
$ws = $wb.Worksheets.Item("Debit-ZeroDollar")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Fri Aug 29 12:39:46 EDT 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Fri Aug 29 12:39:51 EDT 2025"
$ws.Range("A4").Value = "Pass"
$ws.Range("B4").Value = "Fri Aug 29 12:39:56 EDT 2025"
$ws.Range("A5").Value = "Pass"
$ws.Range("B5").Value = "Fri Aug 29 12:40:05 EDT 2025"
$ws.Range("A6").Value = "Pass"
$ws.Range("B6").Value = "Fri Aug 29 12:40:09 EDT 2025"
$ws.Range("A7").Value = "Pass"
$ws.Range("B7").Value = "Fri Aug 29 12:40:14 EDT 2025"
$ws.Range("A8").Value = "Pass"
$ws.Range("B8").Value = "Fri Aug 29 12:40:21 EDT 2025"

# --- Debit-Void ---
$ws = $wb.Worksheets.Item("Debit-Void")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Fri Aug 29 12:40:26 EDT 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Fri Aug 29 12:40:34 EDT 2025"
$ws.Range("A4").Value = "Pass"
$ws.Range("B4").Value = "Fri Aug 29 12:40:43 EDT 2025"
$ws.Range("A5").Value = "Pass"
$ws.Range("B5").Value = "Fri Aug 29 12:40:52 EDT 2025"
$ws.Range("A6").Value = "Pass"
$ws.Range("B6").Value = "Fri Aug 29 12:41:01 EDT 2025"
$ws.Range("A7").Value = "Pass"
$ws.Range("B7").Value = "Fri Aug 29 12:41:10 EDT 2025"
$ws.Range("A8").Value = "Pass"
$ws.Range("B8").Value = "Fri Aug 29 12:41:22 EDT 2025"

# --- Debit-Credit ---
$ws = $wb.Worksheets.Item("Debit-Credit")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Fri Aug 29 12:41:30 EDT 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Fri Aug 29 12:41:39 EDT 2025"
$ws.Range("A4").Value = "Pass"
$ws.Range("B4").Value = "Fri Aug 29 12:41:47 EDT 2025"
$ws.Range("A5").Value = "Pass"
$ws.Range("B5").Value = "Fri Aug 29 12:41:57 EDT 2025"
$ws.Range("A6").Value = "Pass"
$ws.Range("B6").Value = "Fri Aug 29 12:42:10 EDT 2025"
$ws.Range("A7").Value = "Pass"
$ws.Range("B7").Value = "Fri Aug 29 12:42:23 EDT 2025"
$ws.Range("A8").Value = "Pass"
$ws.Range("B8").Value = "Fri Aug 29 12:42:32 EDT 2025"

# --- Debit-Credit-Void ---
$ws = $wb.Worksheets.Item("Debit-Credit-Void")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Fri Aug 29 12:42:42 EDT 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Fri Aug 29 12:42:54 EDT 2025"
$ws.Range("A4").Value = "Pass"
$ws.Range("B4").Value = "Fri Aug 29 12:43:11 EDT 2025"
$ws.Range("A5").Value = "Pass"
$ws.Range("B5").Value = "Fri Aug 29 12:43:25 EDT 2025"
$ws.Range("A6").Value = "Pass"
$ws.Range("B6").Value = "Fri Aug 29 12:43:37 EDT 2025"
$ws.Range("A7").Value = "Pass"
$ws.Range("B7").Value = "Fri Aug 29 12:43:50 EDT 2025"
$ws.Range("A8").Value = "Pass"
$ws.Range("B8").Value = "Fri Aug 29 12:44:03 EDT 2025"

# --- Debit-MRF ---
$ws = $wb.Worksheets.Item("Debit-MRF")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Fri Aug 29 12:44:16 EDT 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Fri Aug 29 12:44:21 EDT 2025"
$ws.Range("A4").Value = "Pass"
$ws.Range("B4").Value = "Fri Aug 29 12:44:27 EDT 2025"
$ws.Range("A5").Value = "Pass"
$ws.Range("B5").Value = "Fri Aug 29 12:44:32 EDT 2025"
$ws.Range("A6").Value = "Pass"
$ws.Range("B6").Value = "Fri Aug 29 12:44:38 EDT 2025"
$ws.Range("A7").Value = "Pass"
$ws.Range("B7").Value = "Fri Aug 29 12:44:43 EDT 2025"
$ws.Range("A8").Value = "Pass"
$ws.Range("B8").Value = "Fri Aug 29 12:44:48 EDT 2025"
$ws.Range("A9").Value = "Pass"
$ws.Range("B9").Value = "Fri Aug 29 12:44:54 EDT 2025"
$ws.Range("A10").Value = "Pass"
$ws.Range("B10").Value = "Fri Aug 29 12:44:59 EDT 2025"
$ws.Range("A11").Value = "Pass"
$ws.Range("B11").Value = "Fri Aug 29 12:45:04 EDT 2025"
$ws.Range("A12").Value = "Pass"
$ws.Range("B12").Value = "Fri Aug 29 12:45:09 EDT 2025"
$ws.Range("A13").Value = "Pass"
$ws.Range("B13").Value = "Fri Aug 29 12:45:15 EDT 2025"
$ws.Range("A14").Value = "Pass"
$ws.Range("B14").Value = "Fri Aug 29 12:45:20 EDT 2025"
$ws.Range("A15").Value = "Pass"
$ws.Range("B15").Value = "Fri Aug 29 12:45:25 EDT 2025"
$ws.Range("A16").Value = "Pass"
$ws.Range("B16").Value = "Fri Aug 29 12:45:31 EDT 2025"
$ws.Range("A17").Value = "Pass"
$ws.Range("B17").Value = "Fri Aug 29 12:45:36 EDT 2025"
$ws.Range("A18").Value = "Pass"
$ws.Range("B18").Value = "Fri Aug 29 12:45:41 EDT 2025"
$ws.Range("A19").Value = "Pass"
$ws.Range("B19").Value = "Fri Aug 29 12:45:46 EDT 2025"
$ws.Range("A20").Value = "Pass"
$ws.Range("B20").Value = "Fri Aug 29 12:45:51 EDT 2025"

# --- Void-MRF ---
$ws = $wb.Worksheets.Item("Void-MRF")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Fri Aug 29 12:45:57 EDT 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Fri Aug 29 12:46:01 EDT 2025"
$ws.Range("A4").Value = "Pass"
$ws.Range("B4").Value = "Fri Aug 29 12:46:05 EDT 2025"
$ws.Range("A5").Value = "Pass"
$ws.Range("B5").Value = "Fri Aug 29 12:46:09 EDT 2025"
$ws.Range("A6").Value = "Pass"
$ws.Range("B6").Value = "Fri Aug 29 12:46:13 EDT 2025"
$ws.Range("A7").Value = "Pass"
$ws.Range("B7").Value = "Fri Aug 29 12:46:17 EDT 2025"
$ws.Range("A8").Value = "Pass"
$ws.Range("B8").Value = "Fri Aug 29 12:46:22 EDT 2025"
$ws.Range("A9").Value = "Pass"
$ws.Range("B9").Value = "Fri Aug 29 12:46:26 EDT 2025"
$ws.Range("A10").Value = "Pass"
$ws.Range("B10").Value = "Fri Aug 29 12:46:30 EDT 2025"
$ws.Range("A11").Value = "Pass"
$ws.Range("B11").Value = "Fri Aug 29 12:46:34 EDT 2025"

# --- Credit-MRF ---
$ws = $wb.Worksheets.Item("Credit-MRF")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Fri Aug 29 12:46:38 EDT 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Fri Aug 29 12:46:42 EDT 2025"
$ws.Range("A4").Value = "Pass"
$ws.Range("B4").Value = "Fri Aug 29 12:46:51 EDT 2025"
$ws.Range("A5").Value = "Pass"
$ws.Range("B5").Value = "Fri Aug 29 12:46:55 EDT 2025"
$ws.Range("A6").Value = "Pass"
$ws.Range("B6").Value = "Fri Aug 29 12:47:00 EDT 2025"
$ws.Range("A7").Value = "Pass"
$ws.Range("B7").Value = "Fri Aug 29 12:47:04 EDT 2025"
$ws.Range("A8").Value = "Pass"
$ws.Range("B8").Value = "Fri Aug 29 12:47:09 EDT 2025"
$ws.Range("A9").Value = "Pass"
$ws.Range("B9").Value = "Fri Aug 29 12:47:14 EDT 2025"
$ws.Range("A10").Value = "Pass"
$ws.Range("B10").Value = "Fri Aug 29 12:47:18 EDT 2025"
$ws.Range("A11").Value = "Pass"
$ws.Range("B11").Value = "Fri Aug 29 12:47:22 EDT 2025"
$ws.Range("A12").Value = "Pass"
$ws.Range("B12").Value = "Fri Aug 29 12:47:27 EDT 2025"
$ws.Range("A13").Value = "Pass"
$ws.Range("B13").Value = "Fri Aug 29 12:47:31 EDT 2025"

# --- DebitCredit-RemID-Pipe ---
$ws = $wb.Worksheets.Item("DebitCredit-RemID-Pipe")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Fri Aug 29 12:47:36 EDT 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Fri Aug 29 12:47:45 EDT 2025"
$ws.Range("A4").Value = "Pass"
$ws.Range("B4").Value = "Fri Aug 29 12:47:53 EDT 2025"
$ws.Range("A5").Value = "Pass"
$ws.Range("B5").Value = "Fri Aug 29 12:48:01 EDT 2025"
$ws.Range("A6").Value = "Pass"
$ws.Range("B6").Value = "Fri Aug 29 12:48:11 EDT 2025"
$ws.Range("A7").Value = "Pass"
$ws.Range("B7").Value = "Fri Aug 29 12:48:19 EDT 2025"
$ws.Range("A8").Value = "Pass"
$ws.Range("B8").Value = "Fri Aug 29 12:48:28 EDT 2025"

# --- Debit-RemID-Pipe ---
$ws = $wb.Worksheets.Item("Debit-RemID-Pipe")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Fri Aug 29 12:48:37 EDT 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Fri Aug 29 12:48:42 EDT 2025"
$ws.Range("A4").Value = "Pass"
$ws.Range("B4").Value = "Fri Aug 29 12:48:47 EDT 2025"
$ws.Range("A5").Value = "Pass"
$ws.Range("B5").Value = "Fri Aug 29 12:48:52 EDT 2025"
$ws.Range("A6").Value = "Pass"
$ws.Range("B6").Value = "Fri Aug 29 12:48:56 EDT 2025"
$ws.Range("A7").Value = "Pass"
$ws.Range("B7").Value = "Fri Aug 29 12:49:01 EDT 2025"
$ws.Range("A8").Value = "Pass"
$ws.Range("B8").Value = "Fri Aug 29 12:49:05 EDT 2025"

# --- DebitVoid-RemID-Pipe ---
$ws = $wb.Worksheets.Item("DebitVoid-RemID-Pipe")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Fri Aug 29 12:49:10 EDT 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Fri Aug 29 12:49:18 EDT 2025"
$ws.Range("A4").Value = "Pass"
$ws.Range("B4").Value = "Fri Aug 29 12:49:27 EDT 2025"
$ws.Range("A5").Value = "Pass"
$ws.Range("B5").Value = "Fri Aug 29 12:49:36 EDT 2025"
$ws.Range("A6").Value = "Pass"
$ws.Range("B6").Value = "Fri Aug 29 12:49:44 EDT 2025"
$ws.Range("A7").Value = "Pass"
$ws.Range("B7").Value = "Fri Aug 29 12:49:53 EDT 2025"
$ws.Range("A8").Value = "Pass"
$ws.Range("B8").Value = "Fri Aug 29 12:50:02 EDT 2025"
